$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.910.32'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.360.32'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.71'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("E9").Value = '  -3.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.124'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0784'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.35'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").Value = '2.725.72'
$ws.Range("E15").Value = '  +1.51%  '
$ws.Range("D16").Value = '2.364.64'
$ws.Range("E16").Value = '  +2.79%  '
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").Value = '42.912.84'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("D21").Value = '0.0₃0883'
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("E24").Value = '  -2.29%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("E28").Value = '  +14.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.27'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.17%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("E34").Value = '  +4.37%  '
$ws.Range("E35").Value = '  +6.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '129.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.42%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.81%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("E40").Value = '  -2.68%  '
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.91'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.25%  '
$ws.Range("D43").Value = '1.926.96'
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("E45").Value = '  +3.19%  '
$ws.Range("E46").Value = '  -8.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '2.587.36'
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("E49").Value = '  +2.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("E51").Value = '  -3.08%  '
